$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.427515666666666
$ws.Range("H2").Value = 19.282547
$ws.Range("I2").Value = 0.399676466497628
$ws.Range("J2").Value = 0.399676466497628
$ws.Range("M2").Value = 43.28265566666667
$ws.Range("N2").Value = 129.847967
$ws.Range("O2").Value = 0.667219228070094
$ws.Range("P2").Value = 0.667219228070094
$ws.Range("Q2").Value = 278.1999473924387
$ws.Range("R2").Value = 2503.799526531949
$ws.Range("S2").Value = 0.2666718234543302
$ws.Range("T2").Value = 0.2666718234543302

$ws.Range("G3").Value = 6.427515666666666
$ws.Range("H3").Value = 19.282547
$ws.Range("I3").Value = 0.399676466497628
$ws.Range("J3").Value = 0.399676466497628
$ws.Range("O3").Value = 0.2872263480299067
$ws.Range("P3").Value = 0.2872263480299067
$ws.Range("Q3").Value = 119.7602700131535
$ws.Range("R3").Value = 1077.842430118382
$ws.Range("S3").Value = 0.114797611865611
$ws.Range("T3").Value = 0.114797611865611

$ws.Range("G4").Value = 6.427515666666666
$ws.Range("H4").Value = 19.282547
$ws.Range("I4").Value = 0.399676466497628
$ws.Range("J4").Value = 0.399676466497628
$ws.Range("M4").Value = 2.955125333333334
$ws.Range("N4").Value = 8.865376000000001
$ws.Range("O4").Value = 0.04555442389999943
$ws.Range("P4").Value = 0.04555442389999944
$ws.Range("Q4").Value = 18.99411437696356
$ws.Range("R4").Value = 170.947029392672
$ws.Range("S4").Value = 0.01820703117768687
$ws.Range("T4").Value = 0.01820703117768687

$ws.Range("I5").Value = 0.1137908927671639
$ws.Range("J5").Value = 0.1137908927671639
$ws.Range("M5").Value = 43.28265566666667
$ws.Range("N5").Value = 129.847967
$ws.Range("O5").Value = 0.667219228070094
$ws.Range("P5").Value = 0.667219228070094
$ws.Range("Q5").Value = 79.20561512908468
$ws.Range("R5").Value = 712.8505361617621
$ws.Range("S5").Value = 0.07592347163351393
$ws.Range("T5").Value = 0.07592347163351393

$ws.Range("I6").Value = 0.1137908927671639
$ws.Range("J6").Value = 0.1137908927671639
$ws.Range("O6").Value = 0.2872263480299067
$ws.Range("P6").Value = 0.2872263480299067
$ws.Range("S6").Value = 0.0326837425685752
$ws.Range("T6").Value = 0.0326837425685752

$ws.Range("I7").Value = 0.1137908927671639
$ws.Range("J7").Value = 0.1137908927671639
$ws.Range("M7").Value = 2.955125333333334
$ws.Range("N7").Value = 8.865376000000001
$ws.Range("O7").Value = 0.04555442389999943
$ws.Range("P7").Value = 0.04555442389999944
$ws.Range("Q7").Value = 5.407767065237334
$ws.Range("R7").Value = 48.66990358713601
$ws.Range("S7").Value = 0.005183678565074763
$ws.Range("T7").Value = 0.005183678565074764

$ws.Range("G8").Value = 7.824318999999999
$ws.Range("H8").Value = 23.472957
$ws.Range("I8").Value = 0.486532640735208
$ws.Range("J8").Value = 0.486532640735208
$ws.Range("M8").Value = 43.28265566666667
$ws.Range("N8").Value = 129.847967
$ws.Range("O8").Value = 0.667219228070094
$ws.Range("P8").Value = 0.667219228070094
$ws.Range("Q8").Value = 338.6573051031577
$ws.Range("R8").Value = 3047.915745928419
$ws.Range("S8").Value = 0.3246239329822499
$ws.Range("T8").Value = 0.3246239329822499

$ws.Range("G9").Value = 7.824318999999999
$ws.Range("H9").Value = 23.472957
$ws.Range("I9").Value = 0.486532640735208
$ws.Range("J9").Value = 0.486532640735208
$ws.Range("O9").Value = 0.2872263480299067
$ws.Range("P9").Value = 0.2872263480299067
$ws.Range("Q9").Value = 145.7861177948713
$ws.Range("R9").Value = 1312.075060153842
$ws.Range("S9").Value = 0.1397449935957204
$ws.Range("T9").Value = 0.1397449935957204

$ws.Range("G10").Value = 7.824318999999999
$ws.Range("H10").Value = 23.472957
$ws.Range("I10").Value = 0.486532640735208
$ws.Range("J10").Value = 0.486532640735208
$ws.Range("M10").Value = 2.955125333333334
$ws.Range("N10").Value = 8.865376000000001
$ws.Range("O10").Value = 0.04555442389999943
$ws.Range("P10").Value = 0.04555442389999944
$ws.Range("Q10").Value = 23.12184329298133
$ws.Range("R10").Value = 208.096589636832
$ws.Range("S10").Value = 0.0221637141572378
$ws.Range("T10").Value = 0.0221637141572378

